# Update the "Förändrad" (changed/updated) date in column C for every
# data row (rows 2 through 481) from 45192 (2023-09-23) to 45202 (2023-10-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 481
$newValue = 45202

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newValue
}
